# Applies the "EC" (Estado de Cuenta) worksheet update:
#  - updates VALOR MORA total, worker count and periods count
#  - replaces the worker detail table (rows 16-21) with the new, larger
#    table (rows 16-27), pushing the signature block down accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary header fields -------------------------------------------------
$ws.Range("E11").Value = 1116158   # VALOR MORA
$ws.Range("C13").Value = 9         # Cant. Trabajadores
$ws.Range("F13").Value = 10        # Cant. Periodos

# ---- Make room for the extra detail rows -----------------------------------
# Original table occupies rows 16-21 (6 rows); new table needs rows 16-27
# (12 rows), so insert 6 new rows right after the last existing data row
# (21), before the first blank row (22). This pushes the signature block
# (previously rows 26-27) down to rows 32-33, exactly like the target file.
$ws.Range("A22:A27").EntireRow.Insert()

# Row 21 currently still carries the special bottom-border styling that used
# to mark the last row of the (shorter) table. Move that styling onto row 27,
# which is now the new last row of the table.
$ws.Range("B21:J21").Copy()
$ws.Range("B27:J27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Copy the formatting of a normal data row (20) onto row 21 and the freshly
# inserted rows (22-26), so they all look like regular table rows.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Fill in the new worker detail table (rows 16-27) ----------------------
$data = @(
    @("CC", "9148943",    "RONALD JAIR MARTINEZ CABARCAS",  "2507", 219810, 5495252),
    @("CC", "9148943",    "RONALD JAIR MARTINEZ CABARCAS",  "2506", 219810, 5495252),
    @("CC", "73112875",   "YIMY FRANCISCO OROZCO CORDOBA",  "2505", 12467,  9350000),
    @("CC", "73112875",   "YIMY FRANCISCO OROZCO CORDOBA",  "2503", 12467,  9350000),
    @("CC", "92226976",   "JULIO JAVIER PUERTA FLOREZ",     "2505", 17272,  18505500),
    @("CC", "73184394",   "JOSE HERNEY RIVERA VERGARA",     "2101", 147333, 8500000),
    @("CC", "8665244",    "ISRAEL SANCHEZ JURADO",          "1909", 160000, 4000000),
    @("CC", "8665244",    "ISRAEL SANCHEZ JURADO",          "1908", 90666,  4000000),
    @("CC", "1096184009", "HECTOR FABIAN DIAZ GARCIA",      "2309", 6867,   6029579),
    @("CC", "1216714256", "JUAN CARLOS CASSIANI GUETTE",    "2208", 122666, 3710000),
    @("CC", "73181412",   "JAIME ALBERTO WATTS ROSSI",      "1909", 51333,  3710000),
    @("CC", "91326446",   "ABEL DEL CARMEN SIDEROL PACHECO","2012", 55467,  2600000)
)

$row = 16
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $entry[1]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $entry[2]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $entry[3]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $entry[4]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value = $entry[5]   # G - Salario Basico
    $row = $row + 1
}
